$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214 (pushing existing rows 214-252 down to 215-253)
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new weekly record
$ws.Cells.Item(214, 1).Value = 4
$ws.Cells.Item(214, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(214, 3).Value = "Los Lagos"
$ws.Cells.Item(214, 4).Value = 44711
$ws.Cells.Item(214, 5).Value = 10
$ws.Cells.Item(214, 6).Value = 100112044
$ws.Cells.Item(214, 7).Value = "Perejil"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 70
$ws.Cells.Item(214, 11).Value = 6000
$ws.Cells.Item(214, 12).Value = 6000
$ws.Cells.Item(214, 13).Value = 6000
$ws.Cells.Item(214, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(214, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(214, 16).Value = 3000
$ws.Cells.Item(214, 17).Value = 2
$ws.Cells.Item(214, 18).Value = "Hortaliza"
